# Auto-generated edit script applying the Diabolos_Profits market-data refresh
# (cached value updates to columns H-N across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H28").Value = 56301.332
$ws.Range("I28").Value = 59555
$ws.Range("K28").Value = 59555
$ws.Range("M28").Value = -59070
$ws.Range("H80").Value = 1430570.5
$ws.Range("J80").Value = 4935.8
$ws.Range("L80").Value = 14807.4
$ws.Range("N80").Value = -16803.4
$ws.Range("H83").Value = 1430570.5
$ws.Range("J83").Value = 4935.8
$ws.Range("L83").Value = 44422.2
$ws.Range("N83").Value = -54406.2
$ws.Range("H86").Value = 8378447
$ws.Range("I86").Value = 2048.3076
$ws.Range("K86").Value = 2048.3076
$ws.Range("M86").Value = -925.3076000000001
$ws.Range("H88").Value = 533
$ws.Range("J88").Value = 362.75
$ws.Range("L88").Value = 362.75
$ws.Range("N88").Value = -1174.75
$ws.Range("H89").Value = 8378447
$ws.Range("I89").Value = 2048.3076
$ws.Range("K89").Value = 10241.538
$ws.Range("M89").Value = -4625.538
$ws.Range("H91").Value = 533
$ws.Range("J91").Value = 362.75
$ws.Range("L91").Value = 362.75
$ws.Range("N91").Value = -3170.75
$ws.Range("H92").Value = 182100.36
$ws.Range("I92").Value = 1589.375
$ws.Range("K92").Value = 1589.375
$ws.Range("M92").Value = -341.375
$ws.Range("H96").Value = 125297.375
$ws.Range("I96").Value = 166987.5
$ws.Range("J96").Value = 227
$ws.Range("K96").Value = 500962.5
$ws.Range("L96").Value = 681
$ws.Range("M96").Value = -499589.5
$ws.Range("N96").Value = -3427
$ws.Range("H100").Value = 7923
$ws.Range("I100").Value = 10899.8
$ws.Range("J100").Value = 2961.6667
$ws.Range("K100").Value = 10899.8
$ws.Range("L100").Value = 2961.6667
$ws.Range("M100").Value = -10358.8
$ws.Range("N100").Value = -4043.6667
$ws.Range("H115").Value = 360.66666
$ws.Range("I115").Value = 266.25
$ws.Range("J115").Value = 549.5
$ws.Range("K115").Value = 798.75
$ws.Range("L115").Value = 1648.5
$ws.Range("M115").Value = 768.25
$ws.Range("N115").Value = -4782.5
$ws.Range("H132").Value = 4026.4407
$ws.Range("I132").Value = 3603.44
$ws.Range("K132").Value = 10810.32
$ws.Range("M132").Value = -8280.32

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H13").Value = 3799.2
$ws.Range("J13").Value = 3666
$ws.Range("L13").Value = 3666
$ws.Range("N13").Value = -3954
$ws.Range("H32").Value = 3087.658
$ws.Range("I32").Value = 3087.658
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3087.658
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2800.658
$ws.Range("N32").ClearContents()
$ws.Range("H132").Value = 2062.1333
$ws.Range("I132").Value = 841.9231
$ws.Range("K132").Value = 2525.7693
$ws.Range("M132").Value = 4.23070000000007

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 7920.2573
$ws.Range("I20").Value = 10496.04
$ws.Range("K20").Value = 10496.04
$ws.Range("M20").Value = -10249.04
$ws.Range("H94").Value = 27781428
$ws.Range("I94").Value = 35717410
$ws.Range("J94").Value = 5500
$ws.Range("K94").Value = 35717410
$ws.Range("L94").Value = 5500
$ws.Range("M94").Value = -35716959
$ws.Range("N94").Value = -6402
$ws.Range("H122").Value = 24999.5
$ws.Range("J122").Value = 24999.5
$ws.Range("L122").Value = 24999.5
$ws.Range("N122").Value = -34799.5
$ws.Range("H134").Value = 2685.9062
$ws.Range("I134").Value = 2514.5
$ws.Range("K134").Value = 7543.5
$ws.Range("M134").Value = -5008.5

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 696.5454999999999
$ws.Range("I16").Value = 658.2857
$ws.Range("K16").Value = 658.2857
$ws.Range("M16").Value = -371.2857
$ws.Range("H88").Value = 17187.5
$ws.Range("J88").Value = 17187.5
$ws.Range("L88").Value = 17187.5
$ws.Range("N88").Value = -17999.5
$ws.Range("H91").Value = 17187.5
$ws.Range("J91").Value = 17187.5
$ws.Range("L91").Value = 17187.5
$ws.Range("N91").Value = -19995.5
$ws.Range("H113").Value = 696.5454999999999
$ws.Range("I113").Value = 658.2857
$ws.Range("K113").Value = 658.2857
$ws.Range("M113").Value = 1511.7143
$ws.Range("H134").Value = 1552.3405
$ws.Range("I134").Value = 994.6667
$ws.Range("K134").Value = 2984.0001
$ws.Range("M134").Value = -449.0001000000002
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H81").Value = 7443.778
$ws.Range("I81").Value = 5399
$ws.Range("J81").Value = 9999.75
$ws.Range("K81").Value = 16197
$ws.Range("L81").Value = 29999.25
$ws.Range("M81").Value = -15074
$ws.Range("N81").Value = -32245.25
$ws.Range("H84").Value = 7443.778
$ws.Range("I84").Value = 5399
$ws.Range("J84").Value = 9999.75
$ws.Range("K84").Value = 48591
$ws.Range("L84").Value = 89997.75
$ws.Range("M84").Value = -42975
$ws.Range("N84").Value = -101229.75

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H113").Value = 2684.1765
$ws.Range("J113").Value = 5000
$ws.Range("L113").Value = 5000
$ws.Range("N113").Value = -9340
$ws.Range("H122").Value = 2425.1875
$ws.Range("I122").Value = 2082.182
$ws.Range("K122").Value = 6246.545999999999
$ws.Range("M122").Value = -3796.545999999999
$ws.Range("H132").Value = 3325.9143
$ws.Range("I132").Value = 3834.261
$ws.Range("K132").Value = 11502.783
$ws.Range("M132").Value = -8972.782999999999
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H22").Value = 887.1667
$ws.Range("I22").Value = 540.6667
$ws.Range("K22").Value = 540.6667
$ws.Range("M22").Value = -245.6667
$ws.Range("H27").Value = 887.1667
$ws.Range("I27").Value = 540.6667
$ws.Range("K27").Value = 540.6667
$ws.Range("M27").Value = -433.6667
$ws.Range("H61").Value = 3501
$ws.Range("I61").Value = 3501
$ws.Range("K61").Value = 3501
$ws.Range("M61").Value = -3299
$ws.Range("H68").Value = 4420
$ws.Range("I68").Value = 2100
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 5000
$ws.Range("M68").Value = -1351
$ws.Range("N68").Value = -6498
$ws.Range("H71").Value = 4420
$ws.Range("I71").Value = 2100
$ws.Range("J71").Value = 5000
$ws.Range("K71").Value = 10500
$ws.Range("L71").Value = 25000
$ws.Range("M71").Value = -6756
$ws.Range("N71").Value = -32488
$ws.Range("H82").Value = 2883.5454
$ws.Range("J82").Value = 2828.3333
$ws.Range("L82").Value = 2828.3333
$ws.Range("N82").Value = -3550.3333
$ws.Range("H85").Value = 2883.5454
$ws.Range("J85").Value = 2828.3333
$ws.Range("L85").Value = 2828.3333
$ws.Range("N85").Value = -5324.3333
$ws.Range("H93").Value = 37041436
$ws.Range("I93").Value = 41670990
$ws.Range("J93").Value = 5000
$ws.Range("K93").Value = 41670990
$ws.Range("L93").Value = 5000
$ws.Range("M93").Value = -41669742
$ws.Range("N93").Value = -7496
$ws.Range("H113").Value = 3501
$ws.Range("I113").Value = 3501
$ws.Range("K113").Value = 3501
$ws.Range("M113").Value = -1331
$ws.Range("H122").Value = 3703.5
$ws.Range("I122").Value = 3204.375
$ws.Range("K122").Value = 9613.125
$ws.Range("M122").Value = -7163.125

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H46").Value = 81171.336
$ws.Range("J46").Value = 81171.336
$ws.Range("L46").Value = 81171.336
$ws.Range("N46").Value = -81633.336
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H68").Value = 44733
$ws.Range("H71").Value = 44733
$ws.Range("H86").Value = 49999.668
$ws.Range("J86").Value = 49999.668
$ws.Range("L86").Value = 49999.668
$ws.Range("N86").Value = -52245.668
$ws.Range("H89").Value = 49999.668
$ws.Range("J89").Value = 49999.668
$ws.Range("L89").Value = 249998.34
$ws.Range("N89").Value = -261230.34
$ws.Range("H107").Value = 921.61536
$ws.Range("I107").Value = 725.0909
$ws.Range("K107").Value = 2175.2727
$ws.Range("M107").Value = -255.2727
$ws.Range("H113").Value = 1392.2142
$ws.Range("I113").Value = 462.7857
$ws.Range("J113").Value = 2321.6428
$ws.Range("K113").Value = 1388.3571
$ws.Range("L113").Value = 6964.928400000001
$ws.Range("M113").Value = 781.6428999999998
$ws.Range("N113").Value = -11304.9284
$ws.Range("H134").Value = 81171.336
$ws.Range("J134").Value = 81171.336
$ws.Range("L134").Value = 243514.008
$ws.Range("N134").Value = -248584.008
